$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) date column from 45184 (2023-09-15) to
#    45186 (2023-09-17) for every data row (rows 2..116).
for ($r = 2; $r -le 116; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# 2) Add a friendly display-text second argument to the HYPERLINK formulas
#    in rows 2-4 for columns S, T, V, W, X, Y. The friendly text is the
#    "Beteckning" value from column A of the same row.
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")
foreach ($r in 2..4) {
    $beteckning = $ws.Range("A$r").Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$r")
        $oldFormula = $cell.Formula
        # oldFormula looks like: =HYPERLINK("https://...")
        # turn into:             =HYPERLINK("https://...", "A 59332-2022")
        # (plain substring edit - drop the trailing ")" and append the
        # friendly-name argument - avoids any -replace templating pitfalls)
        $withoutClosingParen = $oldFormula.Substring(0, $oldFormula.Length - 1)
        $newFormula = $withoutClosingParen + ', "' + $beteckning + '")'
        $cell.Formula = $newFormula
    }
}
